# Update "想去人数" (F column) counters across the three category sheets
# (展览 / 演出 / 本地生活) and their combined roll-up sheet (全部类型), then
# append one brand-new 本地生活 (Local Life) event as row 9.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) ----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(12, 6).Value = 1825
$ws1.Cells.Item(14, 6).Value = 305
$ws1.Cells.Item(17, 6).Value = 6257
$ws1.Cells.Item(18, 6).Value = 239
$ws1.Cells.Item(19, 6).Value = 90
$ws1.Cells.Item(21, 6).Value = 3399
$ws1.Cells.Item(22, 6).Value = 882
$ws1.Cells.Item(32, 6).Value = 1317
$ws1.Cells.Item(33, 6).Value = 807
$ws1.Cells.Item(35, 6).Value = 88
$ws1.Cells.Item(38, 6).Value = 1498

# ---- 演出 (Performance) ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(11, 6).Value = 164
$ws2.Cells.Item(13, 6).Value = 80
$ws2.Cells.Item(16, 6).Value = 148
$ws2.Cells.Item(17, 6).Value = 338

# ---- 本地生活 (Local Life) --------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4, 6).Value = 260

# Append new row 9 - duplicate the row-8 formatting for column A (bold,
# centered, bordered index style) onto row 9, then overwrite the values.
$ws3.Cells.Item(8, 1).Copy($ws3.Cells.Item(9, 1))
$ws3.Cells.Item(9, 1).Value = 8
$ws3.Cells.Item(9, 2).NumberFormat = "@"
$ws3.Cells.Item(9, 2).Value = "2024-11-02"
$ws3.Cells.Item(9, 2).Style = "Normal"
$ws3.Cells.Item(9, 3).Value = "北京·桑语映画COSPLAY自拍馆11月2日·一日店长·王小泣"
$ws3.Cells.Item(9, 4).Value = "崇文门外大街3~5号 北京新世界百货崇文门店"
$ws3.Cells.Item(9, 5).Value = "2024.11.02 13:00-11.02 18:00"
$ws3.Cells.Item(9, 6).Value = 0
$ws3.Cells.Item(9, 7).Value = 288
$ws3.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93372"
$ws3.Cells.Item(9, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/AHFSIla51728714804851.jpeg"

# ---- 全部类型 (All types roll-up) ------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 260
$ws4.Cells.Item(21, 6).Value = 1825
$ws4.Cells.Item(22, 6).Value = 164
$ws4.Cells.Item(24, 6).Value = 305
$ws4.Cells.Item(26, 6).Value = 6257
$ws4.Cells.Item(27, 6).Value = 239
$ws4.Cells.Item(28, 6).Value = 90
$ws4.Cells.Item(30, 6).Value = 3399
$ws4.Cells.Item(31, 6).Value = 882
$ws4.Cells.Item(39, 6).Value = 1317
$ws4.Cells.Item(40, 6).Value = 338
$ws4.Cells.Item(43, 6).Value = 807
$ws4.Cells.Item(45, 6).Value = 88
